$d = $word.ActiveDocument

# 1. Title: "You’re invited to our Deriv Partner Seminar"
$d.Content.Find.Execute(
    "You’re invited to our Deriv Partner Seminar", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Zapraszamy Państwa na nasze seminarium dla partnerów Deriv", 2)

# 2. "Dear [PARTNER NAME], "
$d.Content.Find.Execute(
    "Dear [PARTNER NAME], ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Szanowni Państwo [PARTNER NAME], ", 2)

# 3. "We’re excited to let you know that the Deriv Affiliate team will be in [CITY] in [MONTH] to meet with you, our valued partners!"
$d.Content.Find.Execute(
    "We’re excited to let you know that the Deriv Affiliate team will be in [CITY] in [MONTH] to meet with you, our valued partners!",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Z radością informujemy, że zespół Deriv Affiliate będzie obecny w mieście: [CITY] w miesiącu: [MONTH] , aby spotkać się z naszymi cenionymi partnerami!", 2)

# 4. "In this one-day seminar, ..."
$d.Content.Find.Execute(
    "In this one-day seminar, we’ll be providing technical and marketing support, offering the opportunity to network with other partners over a delicious lunch as well as listening to your feedback about our partnership programmes. This is your chance to get your voice heard, which will help us plan future efforts to support you better. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Podczas tego jednodniowego seminarium zapewnimy wsparcie techniczne i marketingowe, zaoferujemy możliwość nawiązania kontaktów z innymi partnerami podczas pysznego lunchu, a także wysłuchamy Państwa opinii na temat naszych programów partnerskich. Jest to szansa, aby Państwa głos został wysłuchany, co pomoże nam zaplanować przyszłe działania, aby lepiej Państwa wspierać. ", 2)

# 5. ". Please note that attendance is confirmed on a first come, first served basis. We look forward to seeing you there!"
$d.Content.Find.Execute(
    ". Please note that attendance is confirmed on a first come, first served basis. We look forward to seeing you there!",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ". Please note that attendance is confirmed on a first come, first served basis. Z niecierpliwością czekamy na Państwa!", 2)
